# Weekly update: insert 4 new daily price rows for Cilantro (Mercado Mayorista
# Lo Valledor de Santiago) right above the existing row 1082, shifting the
# existing data (old rows 1082-1126) down to rows 1086-1130.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows before the current row 1082 (shifts 1082:1126 -> 1086:1130)
$ws.Range("A1082:A1085").EntireRow.Insert()

# Constant columns shared by every data row in this sheet
$marketId = 6
$marketName = "Mercado Mayorista Lo Valledor de Santiago"
$region = "Metropolitana"
$codreg = 13
$categoryId = 100112040
$category = "Cilantro"
$variety = "Sin especificar"
$clasificacion = "Hortaliza"

# row, fecha, calidad, volumen, precioMin, precioMax, precioProm, unidad, origen, precioKg, kgUnidades
$newRows = @(
    @(1082, 44939, "Primera", 250, 9000,  9000,  9000,  "$/caja 36 atados",   "Provincia de Quillota",  250,  36),
    @(1083, 44939, "Primera", 680, 7500,  8000,  7728,  "$/caja 36 atados",   "Región Metropolitana",   215,  36),
    @(1084, 44939, "Primera", 560, 17000, 18000, 17393, "$/docena de atados", "Región Metropolitana",   5798, 3),
    @(1085, 44939, "Primera", 230, 22000, 22000, 22000, "$/docena de atados", "Región de Valparaíso",   7333, 3)
)

foreach ($r in $newRows) {
    $row = $r[0]
    $ws.Cells.Item($row, 1).Value = $marketId
    $ws.Cells.Item($row, 2).Value = $marketName
    $ws.Cells.Item($row, 3).Value = $region
    $ws.Cells.Item($row, 4).Value = $r[1]
    $ws.Cells.Item($row, 5).Value = $codreg
    $ws.Cells.Item($row, 6).Value = $categoryId
    $ws.Cells.Item($row, 7).Value = $category
    $ws.Cells.Item($row, 8).Value = $variety
    $ws.Cells.Item($row, 9).Value = $r[2]
    $ws.Cells.Item($row, 10).Value = $r[3]
    $ws.Cells.Item($row, 11).Value = $r[4]
    $ws.Cells.Item($row, 12).Value = $r[5]
    $ws.Cells.Item($row, 13).Value = $r[6]
    $ws.Cells.Item($row, 14).Value = $r[7]
    $ws.Cells.Item($row, 15).Value = $r[8]
    $ws.Cells.Item($row, 16).Value = $r[9]
    $ws.Cells.Item($row, 17).Value = $r[10]
    $ws.Cells.Item($row, 18).Value = $clasificacion
}
